$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 86, shifting existing rows 86..158 down to 87..159.
$ws.Rows("86").Insert()

# Populate the newly inserted row 86 with its data.
$ws.Range("A86").Value = 8
$ws.Range("B86").Value = "Terminal La Palmera de La Serena"
$ws.Range("C86").Value = "Coquimbo"
$ws.Range("D86").Value = 45072
$ws.Range("E86").Value = 4
$ws.Range("F86").Value = "Fruta"
$ws.Range("G86").Value = 100109
$ws.Range("H86").Value = "Uva"
$ws.Range("I86").Value = 100109001
$ws.Range("J86").Value = "Uva"
$ws.Range("K86").Value = "Red Globe"
$ws.Range("L86").Value = "Primera"
$ws.Range("M86").Value = 500
$ws.Range("N86").Value = 8000
$ws.Range("O86").Value = 9000
$ws.Range("P86").Value = 8500
$ws.Range("Q86").Value = "`$/bandeja 18 kilos"
$ws.Range("R86").Value = "Provincia del Elquí"
$ws.Range("S86").Value = 472
$ws.Range("T86").Value = 18
